$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# EPBDS-8724: Simple Rules / Simple Lookups must not support StringRanges
# for backward compatibility. The test workbook used columns K:L as a
# "StringRange" variant of the B:C (and G:H) sample data for rows
# 164-201; that StringRange block is removed entirely (rows 179-183 in
# that range never had K:L content to begin with).
#
# Two of those rows (164 and 184) had their K:L cells merged, so unmerge
# first and then clear contents+formatting so the <c> nodes disappear
# from the saved XML rather than just losing their value.

$ws.Range("K164:L164").UnMerge()
$ws.Range("K184:L184").UnMerge()

$ws.Range("K164:L178").Clear()
$ws.Range("K184:L201").Clear()
